$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dates (A2/A3): 2025-12-01 -> 2025-12-03 ---
# Use a leading apostrophe so Excel stores it as literal text (not an
# auto-converted date serial), then strip the resulting quote-prefix
# style back off so the cell's style index is unchanged (style 0).
$ws.Range("A2").Value = "'2025-12-03"
$ws.Range("A2").Style = "Normal"

$ws.Range("A3").Value = "'2025-12-03"
$ws.Range("A3").Style = "Normal"

# --- Row 2 (Oklo Inc. / OKLO) numeric updates ---
$ws.Range("D2").Value = 92.19
$ws.Range("E2").Value = 40.7
$ws.Range("F2").Value = 2.95
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 56.6
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 65.32892478746797
$ws.Range("O2").Value = "🟢 상승 우위 (다소 완화)"

# --- Row 3 (NuScale Power Corporation / SMR) numeric updates ---
$ws.Range("D3").Value = 18.94
$ws.Range("E3").Value = 25.2
$ws.Range("F3").Value = -5.02
$ws.Range("H3").Value = 70
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 49.6
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 65.32892478746797
$ws.Range("O3").Value = "🟢 상승 우위 (다소 완화)"
